# Punto de partida, Se agrego escenarios de diferentes precios
#
# This script reproduces (on the "HighOil_AEO_FINAL" sheet):
#  - P8's formula switched from "=+R42" (a flat USD2019/gal AVGAS price) to
#    "=+R43" (that price converted to USD2019/MillonBTU), which cascades
#    through the whole P8:P38 / P7 price-scenario column.
#  - A new "AVGAS" labelled column (R6 header, S7:S38 pasted-as-values
#    mirror of the recomputed P column) added to the sheet.
#  - The previously-active sheet tab moves from "HighOil_AEO" to
#    "HighOil_AEO_FINAL".

$wb = $excel.ActiveWorkbook

$wsFinal = $wb.Worksheets.Item("HighOil_AEO_FINAL")

# --- Core formula edit: point the AVGAS anchor cell at the converted
#     (USD2019/MillonBTU) reference value instead of the raw USD2019/gal one.
$wsFinal.Range("P8").Formula = "=+R43"

# --- New "AVGAS" column: header label in R6 ...
$wsFinal.Range("R6").Value = "AVGAS"

# ... and S7:S38 filled with the (now recalculated) P-column scenario
# prices, pasted in as plain values.
for ($r = 7; $r -le 38; $r++) {
    $cell = $wsFinal.Cells.Item($r, 16)   # column P
    $wsFinal.Cells.Item($r, 19).Value = $cell.Value2   # column S
}

# --- Selection / active tab housekeeping: the workbook now opens on
#     "HighOil_AEO_FINAL" with P38:S38 highlighted, instead of
#     "HighOil_AEO".
$wsFinal.Activate() | Out-Null
$wsFinal.Range("P38:S38").Select() | Out-Null
